$wb = $excel.ActiveWorkbook

# The CI "Generate Report" run moved this file from "Ready for handoff"
# to "In Translation" -- update the Status value everywhere it appears:
# the Overview rollup sheet (columns E/F, one per language) and each
# per-language sheet's "Status" column (column C).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# The new status text is shorter, so those columns were narrowed to fit.
# Column widths are stored as an integer pixel count (5px padding, 6px per
# character at this workbook's font metrics), so ColumnWidth values snap to
# the nearest 1/6th -- 12.5 is the input that lands closest to the
# canonical target width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
